# Weekly update: insert the newest price record for "Vega Modelo de Temuco -
# Locoto" at the top of the data series (row 60). Every existing record
# (rows 60-101) shifts down by one row, which also moves the previously
# last record into a brand-new row 102 and grows the used range to A1:R102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 60..101 down to 61..102, opening up a blank row 60.
$ws.Rows(60).Insert()

# Fill the newly opened row 60 with this week's record. The non-varying,
# series-level columns (Mercado ID, Mercado, Region, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Unidad de comercializacion, Origen, Kg o
# Unidades, Clasificacion) repeat the same values as every other row in
# this series; only the date (D), volume (J) and prices (K/L/M/P) are new.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 45176
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100112042
$ws.Range("G60").Value = "Locoto"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 90
$ws.Range("K60").Value = 2200
$ws.Range("L60").Value = 2200
$ws.Range("M60").Value = 2200
$ws.Range("N60").Value = "$/kilo"
$ws.Range("O60").Value = "Región de Arica y Parinacota"
$ws.Range("P60").Value = 2200
$ws.Range("Q60").Value = 1
$ws.Range("R60").Value = "Hortaliza"
